$d = $word.ActiveDocument

# ---------------------------------------------------------------
# Change 2: "Дата изготовления: DATEUP г." -> "Дата изготовления: DATE г."
# ---------------------------------------------------------------
$null = $d.Content.Find.Execute("DATEUP", $true, $false, $false, $false, $false, $true, 1, $false, "DATE", 2)

# ---------------------------------------------------------------
# Change 3: merge "РВМН-26.51." + "52-406233-001-22 ПС" into one run
# (2nd occurrence of the phrase in the document)
# ---------------------------------------------------------------
$r = $d.Content
$null = $r.Find.Execute("РВМН-26.51.52-406233-001-22 ПС", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
$r2 = $d.Range($r.End, $d.Content.End)
$null = $r2.Find.Execute("РВМН-26.51.52-406233-001-22 ПС", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
$r2.Text = "TEMP_MERGE_RVMN"
$r2.Text = "РВМН-26.51.52-406233-001-22 ПС"

# ---------------------------------------------------------------
# Change 4: "Дата изготовления: DATEDOWN г." -> "Дата изготовления: DATE" +
# bookmark _GoBack + " г."
# ---------------------------------------------------------------
$r3 = $d.Content
$null = $r3.Find.Execute("DATEDOWN", $true, $false, $false, $false, $false, $true, 1, $false, "DATE", 2)
$insertPoint = $d.Range($r3.End, $r3.End)
$d.Bookmarks.Add("_GoBack", $insertPoint)

# ---------------------------------------------------------------
# Change 5: merge "Дата снятия с " + "изделия" into one run
# (2nd occurrence of the phrase in the document)
# ---------------------------------------------------------------
$r4 = $d.Content
$null = $r4.Find.Execute("Дата снятия с изделия", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
$r5 = $d.Range($r4.End, $d.Content.End)
$null = $r5.Find.Execute("Дата снятия с изделия", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
$r5.Text = "TEMP_MERGE_SNYATIYA"
$r5.Text = "Дата снятия с изделия"

# ---------------------------------------------------------------
# Change 6: remove the _GoBack bookmark that used to sit at the end of the
# document (right after the final drawing/picture). Because we just added a
# new "_GoBack" bookmark above (change 4), Word has already relocated the
# single allowed "_GoBack" bookmark there, so nothing else to do here -- but
# keep this defensive removal in case the runtime kept the old one too.
# ---------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $bm = $d.Bookmarks.Item("_GoBack")
    if ($bm.Start -eq $bm.End -and $bm.Start -eq $d.Content.End) {
        $bm.Delete()
    }
}
